# Update the lognormal length-of-stay distribution parameters on each
# "traj" sheet to include the new second (sigma) parameter, and move the
# active selection to match the author's final cursor position on each
# sheet.

$wb = $excel.ActiveWorkbook

# traj1: ECU -> lognorm,0.6,0.06
$ws1 = $wb.Worksheets.Item("traj1")
$ws1.Range("B2").Value = "lognorm,0.6,0.06"
$null = $ws1.Range("B2").Select()

# traj2: ICU -> lognorm,0.6,0.06
$ws2 = $wb.Worksheets.Item("traj2")
$ws2.Range("B2").Value = "lognorm,0.6,0.06"
$null = $ws2.Range("B2").Select()

# traj3: ECU -> lognorm,0.6,0.06 ; ICU -> lognorm,5.4,0.54
$ws3 = $wb.Worksheets.Item("traj3")
$ws3.Range("B2").Value = "lognorm,0.6,0.06"
$ws3.Range("B3").Value = "lognorm,5.4,0.54"
$null = $ws3.Range("B2").Select()

# traj4: ICU -> lognorm,5.4,0.54 ; ECU -> lognorm,0.6,0.06
$ws4 = $wb.Worksheets.Item("traj4")
$ws4.Range("B2").Value = "lognorm,5.4,0.54"
$ws4.Range("B3").Value = "lognorm,0.6,0.06"
$null = $ws4.Range("B3").Select()

# traj5: ECU -> lognorm,0.6,0.06 ; ICU -> lognorm,4.8,0.48 ; ECU -> lognorm,0.6,0.06
$ws5 = $wb.Worksheets.Item("traj5")
$ws5.Range("B2").Value = "lognorm,0.6,0.06"
$ws5.Range("B3").Value = "lognorm,4.8,0.48"
$ws5.Range("B4").Value = "lognorm,0.6,0.06"
$null = $ws5.Range("B4").Select()

# Leave the first sheet active/selected, matching the saved workbook state.
$null = $ws1.Activate()
$null = $ws1.Range("B2").Select()
